# Update the ISL NO (column C), Item Name (column D) and UOM (column E)
# values on Sheet1 to reflect the re-ordered/re-numbered product list.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Each entry: Row, ISL NO, Item Name, UOM
$rows = @(
    @(2,  1, "Desodin 60ml Syrup",                 "60 ml"),
    @(3,  2, "Dinafex 180mg Tablet",                "30's"),
    @(4,  3, "Dinafex 120mg Tablet",                "30's"),
    @(5,  4, "Dinafex 60mg Tablet",                 "30's"),
    @(6,  5, "Dorenta 50mg Tablet",                 "50's"),
    @(7,  6, "Etorix 90mg Tablet",                  "30's"),
    @(8,  7, "Etorix 120mg Tablet",                 "20's"),
    @(9,  8, "Etorix 60mg Tablet - 40's",           "40's"),
    @(10, 9, "Fenobac 100ml Syrup",                 "100ml"),
    @(11, 10, "Flucloxin 500mg Capsule",            "30 's"),
    @(12, 11, "Flucloxin 500mg Capsule - 36's",     "36 's"),
    @(13, 12, "Geminox 320mg Tablet - 8's",         "8 's"),
    @(14, 13, "Ketonic 30mg IM/IV Injection - 4's", "4's"),
    @(15, 14, "Ketonic 10mg Tablet",                "20's"),
    @(16, 15, "Ketonic 30mg Injection",             "5 's"),
    @(17, 16, "Kynol D 25mg Tablet",                "60 's"),
    @(18, 17, "Kynol TR 100mg Capsule",             "50 's"),
    @(19, 18, "Kynol TR 200mg Capsule",             "30 's"),
    @(20, 19, "Naprox Plus 500mg Tablet - 30's",    "30 's"),
    @(21, 20, "Oradin Plus Tablet - 40's",          "40 's"),
    @(22, 21, "Osticare Tablet 24's",                "24's"),
    @(23, 22, "Quinox 100ml IV Infusion",           "1's"),
    @(24, 23, "Rupaday Oral Solution 60ml",         "1's"),
    @(25, 24, "Zithrox 30ml Dry Suspension",        "30ml"),
    @(26, 25, "Zithrox 15ml Suspension",            "15 ml"),
    @(27, 26, "Zithrox 250mg Tablet - 6's",         "6's"),
    @(28, 27, "Zithrox 500mg Tablet",               "6 's")
)

foreach ($entry in $rows) {
    $r = $entry[0]
    $ws.Cells.Item($r, 3).Value = $entry[1]
    $ws.Cells.Item($r, 4).Value = $entry[2]
    $ws.Cells.Item($r, 5).Value = $entry[3]
}
